# Generate Report for Handback
# Applies the "handback" results to the localization-status report:
#  - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#  - Each language sheet's per-file row gets its Latest Target File (linked to
#    the doc), Latest Handback File (xlf name) and Latest Handback DateTime
#    filled in.
#  - A few columns are widened so the newly-populated values are readable.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: zh-cn / de-de status columns
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

# Widen the zh-cn / de-de columns to fit the longer status text.
$overview.Columns.Item(5).ColumnWidth = 29.16666667
$overview.Columns.Item(6).ColumnWidth = 29.16666667

# ---------------------------------------------------------------------
# Helper data for the two language detail sheets
# ---------------------------------------------------------------------
$docMdUrl1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/996a98e1b3a0d8b18348baa31ceed16e811325f2/e2e/ade785ed-c819-45ea-aa12-0d44af389e51.md"
$docMdName1 = "ade785ed-c819-45ea-aa12-0d44af389e51.md"
$docMdUrl2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/996a98e1b3a0d8b18348baa31ceed16e811325f2/e2e/eafc6ada-4ea2-44e7-97ff-d678a18abeef.md"
$docMdName2 = "eafc6ada-4ea2-44e7-97ff-d678a18abeef.md"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("C2").Value = $newStatus
$zh.Range("C3").Value = $newStatus

$zh.Range("J2").Value = "ade785ed-c819-45ea-aa12-0d44af389e51.887f166bf2852eacbe745c85a92a7cfa9f091368.zh-cn.xlf"
$zh.Range("K2").Value = "2016-09-06 18:39:48"

$zh.Range("J3").Value = "eafc6ada-4ea2-44e7-97ff-d678a18abeef.ded132596f11d77ac8f90f07cc1c26f3d8120045.zh-cn.xlf"
$zh.Range("K3").Value = "2016-09-06 18:39:48"

# Rebuild the hyperlinks so the new "Latest Target File" links come out in
# document order (A2, I2, A3, I3), same as the source doc links reused for
# the now-confirmed target files.
$zh.Range("A2").Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Cells.Item(2, 1), $docMdUrl1, "", "", $docMdName1)
$zh.Hyperlinks.Add($zh.Cells.Item(2, 9), $docMdUrl1, "", "", $docMdName1)
$zh.Hyperlinks.Add($zh.Cells.Item(3, 1), $docMdUrl2, "", "", $docMdName2)
$zh.Hyperlinks.Add($zh.Cells.Item(3, 9), $docMdUrl2, "", "", $docMdName2)

$zh.Columns.Item(3).ColumnWidth = 29.16666667
$zh.Columns.Item(9).ColumnWidth = 39.16666667
$zh.Columns.Item(10).ColumnWidth = 39.16666667

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("C2").Value = $newStatus
$de.Range("C3").Value = $newStatus

$de.Range("J2").Value = "ade785ed-c819-45ea-aa12-0d44af389e51.887f166bf2852eacbe745c85a92a7cfa9f091368.de-de.xlf"
$de.Range("K2").Value = "2016-09-06 18:39:56"

$de.Range("J3").Value = "eafc6ada-4ea2-44e7-97ff-d678a18abeef.ded132596f11d77ac8f90f07cc1c26f3d8120045.de-de.xlf"
$de.Range("K3").Value = "2016-09-06 18:39:56"

$de.Range("A2").Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Cells.Item(2, 1), $docMdUrl1, "", "", $docMdName1)
$de.Hyperlinks.Add($de.Cells.Item(2, 9), $docMdUrl1, "", "", $docMdName1)
$de.Hyperlinks.Add($de.Cells.Item(3, 1), $docMdUrl2, "", "", $docMdName2)
$de.Hyperlinks.Add($de.Cells.Item(3, 9), $docMdUrl2, "", "", $docMdName2)

$de.Columns.Item(3).ColumnWidth = 29.16666667
$de.Columns.Item(9).ColumnWidth = 39.16666667
$de.Columns.Item(10).ColumnWidth = 39.16666667
